$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.773.04'
$ws.Range('E2').Value = '  +2.20%  '

$ws.Range('D3').Value = '2.645.67'

$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.91'
$ws.Range('E5').Value = '  +1.82%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.67'
$ws.Range('E6').Value = '  +3.92%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  -0.66%  '

$ws.Range('D9').Value = '2.643.54'
$ws.Range('E9').Value = '  +1.35%  '

$ws.Range('E10').Value = '  +13.19%  '

$ws.Range('E11').Value = '  -0.35%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.23'
$ws.Range('E12').Value = '  +1.02%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.351'
$ws.Range('E13').Value = '  +2.01%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.94'
$ws.Range('E14').Value = '  +2.35%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000189'
$ws.Range('E15').Value = '  +6.30%  '

$ws.Range('D16').Value = '3.132.31'
$ws.Range('E16').Value = '  +1.68%  '

$ws.Range('D17').Value = '68.681.56'
$ws.Range('E17').Value = '  +2.25%  '

$ws.Range('D18').Value = '2.647.25'
$ws.Range('E18').Value = '  +1.56%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.38'
$ws.Range('E19').Value = '  +3.46%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '365.33'
$ws.Range('E20').Value = '  -0.28%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.47'
$ws.Range('E21').Value = '  +1.78%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.28'
$ws.Range('E22').Value = '  -0.16%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.92'
$ws.Range('E23').Value = '  +2.14%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.11'
$ws.Range('E24').Value = '  +4.67%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.87'
$ws.Range('E25').Value = '  +8.01%  '

$ws.Range('E26').Value = '  +0.13%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.09'
$ws.Range('E27').Value = '  +1.50%  '

$ws.Range('E28').Value = '  +7.68%  '

$ws.Range('D29').Value = '2.779.65'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '582.60'
$ws.Range('E30').Value = '  +0.55%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.16%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.02'
$ws.Range('E32').Value = '  +5.09%  '

$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.42'
$ws.Range('E33').Value = '  +3.79%  '

$ws.Range('E34').Value = '  +3.24%  '

$ws.Range('E35').Value = '  +6.02%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.09%  '

$ws.Range('E37').Value = '  +4.55%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '158.96'
$ws.Range('E38').Value = '  +2.04%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.94'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.32'
$ws.Range('E40').Value = '  +1.92%  '

$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.368'
$ws.Range('E41').Value = '  +0.97%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.39'
$ws.Range('E42').Value = '  +3.58%  '

$ws.Range('E43').Value = '  +6.59%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.70'
$ws.Range('E44').Value = '  +5.40%  '

$ws.Range('D45').Value = '0.0₆0321'
$ws.Range('E45').Value = '  +11.15%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.71'
$ws.Range('E46').Value = '  -0.16%  '

$ws.Range('E47').Value = '  +0.10%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '155.95'
$ws.Range('E48').Value = '  +0.93%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.73'
$ws.Range('E49').Value = '  +0.52%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.06'
$ws.Range('E50').Value = '  +3.32%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.71'
$ws.Range('E51').Value = '  +1.42%  '
